# Update cryptos list spreadsheet with latest price/volume data
# (commit: "Updated cryptos list on Tue Sep 26 03:20:59 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.357.23'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.592.50'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.50'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").Value = '1.816.46'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '1.571.51'
$ws.Range("E14").Value = '  -0.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.70'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '26.350.28'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  +4.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.93'
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.09'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.07'
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("D34").Value = '1.334.99'
$ws.Range("E34").Value = '  +4.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("E36").Value = '  -2.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.49'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("E40").Value = '  +5.26%  '
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.994'
$ws.Range("E42").Value = '  -24.25%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.765'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").Value = '1.729.62'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.92'
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.35'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").Value = '  -4.05%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0504'
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0981'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.43'
$ws.Range("E51").Value = '  -0.69%  '
